$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewJoiners")

$passwords = @(
    "PHBIwc4+^Kp9",
    ">m<Ot42gtk`$M",
    "2A#(%46dc4y)",
    "i26G1ji51*cl",
    "i12#!-(()+<4",
    "NJ3909O+*6<P",
    "^3^<sCo2&#HJ",
    "IEa+%tB@(J23",
    "3(^@K%#GdgP&",
    "B^+8`$q%+d#J%",
    "5f`$dB8!cnJE>",
    "cEF4NN+^52d3"
)

for ($i = 0; $i -lt $passwords.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $passwords[$i]
}
